$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3633115.79
$ws.Range("C7").Value = -18.22975667148522
$ws.Range("D7").Value = 3183
$ws.Range("E7").Value = 3183
$ws.Range("F7").Value = 1141.412437951618
$ws.Range("G7").Value = 21.6663124108847
